$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "\31 52170-case-636"
$ws.Range("B4").Value = "\31 52171-case-641"
$ws.Range("C4").Value = "1 TB"
$ws.Range("D4").Value = "Black"
